# Solutions slide (slide 4): expand the first bullet so it reads
# "Create SAI with pyinstaller for multiplatform" instead of
# "Create SAI with pyinstaller".
#
# The run " pyinstaller" (a space followed by "pyinstaller") is replaced
# by " pyinstaller for multiplatform" in a single targeted edit so the
# surrounding runs ("Create", " SAI ", "with") are left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$target = $tr.Characters(16, 12)
$target.Text = " pyinstaller for multiplatform"
